$d = $word.ActiveDocument

$d.Content.Find.Execute("75÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "40÷5=", 2) | Out-Null
$d.Content.Find.Execute("37÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "62÷4=", 2) | Out-Null
$d.Content.Find.Execute("71÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "55÷4=", 2) | Out-Null
$d.Content.Find.Execute("60÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "18÷4=", 2) | Out-Null
$d.Content.Find.Execute("85÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "15÷2=", 2) | Out-Null
$d.Content.Find.Execute("29÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "72÷5=", 2) | Out-Null
$d.Content.Find.Execute("53÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "36÷2=", 2) | Out-Null
$d.Content.Find.Execute("58÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "11÷3=", 2) | Out-Null
$d.Content.Find.Execute("69÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "48÷3=", 2) | Out-Null
$d.Content.Find.Execute("38÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "29÷8=", 2) | Out-Null
$d.Content.Find.Execute("73÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "85÷5=", 2) | Out-Null
$d.Content.Find.Execute("68÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "98÷6=", 2) | Out-Null
$d.Content.Find.Execute("33÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "91÷3=", 2) | Out-Null
$d.Content.Find.Execute("55÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "29÷2=", 2) | Out-Null
$d.Content.Find.Execute("99÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "35÷2=", 2) | Out-Null
$d.Content.Find.Execute("95÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "30÷8=", 2) | Out-Null
$d.Content.Find.Execute("92÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "58÷5=", 2) | Out-Null
$d.Content.Find.Execute("51÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "84÷5=", 2) | Out-Null
$d.Content.Find.Execute("50÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "19÷7=", 2) | Out-Null
$d.Content.Find.Execute("59÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "29÷7=", 2) | Out-Null
$d.Content.Find.Execute("48÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "73÷4=", 2) | Out-Null
$d.Content.Find.Execute("73÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "52÷3=", 2) | Out-Null
$d.Content.Find.Execute("95÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "81÷6=", 2) | Out-Null
$d.Content.Find.Execute("64÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "12÷2=", 2) | Out-Null
$d.Content.Find.Execute("33÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "59÷2=", 2) | Out-Null
